# Adding Array & LinkedList Module
# -----------------------------------------------------------------
# This script recreates, via Excel COM automation, the edit that:
#  - adds a new worksheet "codeInvalid" after "pythonCode"
#  - reshapes sheet "pythonCode": removes old A5 "pythonCode_invalid" cell,
#    adds two more (empty, wrap-styled) rows below the code sample, and
#    widens column A
#  - populates the new "codeInvalid" sheet with an invalid-code sample
#  - leaves "codeInvalid" as the active/selected sheet

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- 1. Reshape "pythonCode" (sheet1) ----------------------------------
# Drop the old row 5 ("pythonCode_invalid") entirely - that text moves to
# the new sheet under a different label.
$ws1.Rows.Item(5).Delete()

# Style rows 3 & 4 col A the same way as A2 (wrap-text style) and widen
# column A to fit the sample code.
$ws1.Range("A3").WrapText = $true
$ws1.Range("A4").WrapText = $true

# Re-autofit row 2 now that only two helper rows remain below it.
$ws1.Rows.Item(2).RowHeight = 46.5

# ColumnWidth adds the standard padding on save, so back the requested
# displayed width (23.5 chars) off by that padding.
$ws1.Columns.Item(1).ColumnWidth = 22.666666666666668

$null = $ws1.Range("A3:B4").Select()

# ---- 2. Add the new "codeInvalid" worksheet ----------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "codeInvalid"

$ws2.Range("A1").Value = "inValid_pythonCode"
$ws2.Range("B1").Value = "Column2"
$ws2.Range("A2").Value = "color = red blue yellow print colors"
$ws2.Range("B2").Value = "Test"

$ws2.Range("A1").WrapText = $true
$ws2.Range("A2").WrapText = $true

$ws2.Rows.Item(1).RowHeight = 46.5
$ws2.Rows.Item(2).RowHeight = 77.5

$null = $ws2.Range("E2").Select()

# ---- 3. Make the new sheet the active tab ------------------------------
$null = $ws2.Activate()
